$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.762.70'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '3.524.07'
$ws.Range('E3').Value = '  +1.25%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.997'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.46'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '196.59'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.93%  '
$ws.Range('E7').Value = '  +0.70%  '
$ws.Range('E8').Value = '  -0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.201'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.650'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '53.85'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000303'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.53'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '4.080.10'
$ws.Range('E14').Value = '  +1.23%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '597.03'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.56%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '12.83'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.91%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '19.13'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.64%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '69.906.57'
$ws.Range('E18').Value = '  +0.74%  '
$ws.Range('D19').Value = '3.518.37'
$ws.Range('E19').Value = '  +1.19%  '
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.994'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '18.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.54%  '
$ws.Range('E23').Value = '  +5.19%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.68'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.82%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '102.45'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.32%  '
$ws.Range('E26').Value = '  +5.70%  '
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.62'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.87%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.51'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('B30').Value = 'NEARProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.08'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.76%  '
$ws.Range('B31').Value = 'dogwifhat'
$ws.Range('C31').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.30'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.116'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.32%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '63.12'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.05%  '
$ws.Range('D35').Value = '0.0₃0858'
$ws.Range('E35').Value = '  +10.96%  '
$ws.Range('D36').Value = '3.738.54'
$ws.Range('E36').Value = '  +3.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.07'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.42%  '
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('E39').Value = '  +0.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.393'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.63'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '487.59'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.89%  '
$ws.Range('E43').Value = '  -3.48%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0455'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.36%  '
$ws.Range('E45').Value = '  -1.13%  '
$ws.Range('E46').Value = '  -3.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.29'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.49%  '
$ws.Range('E48').Value = '  +0.21%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.52'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.24%  '
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('E51').Value = '  +11.54%  '
